# ------------------------------------------------------------------
# Program_Health_Report.pptx edit script
# 1) Bump the AI-generated date stamp on the title slide subtitle.
# 2) Replace the Executive Program Status AI narrative (slide 2) with
#    the new "Immediate Improvement Areas" analysis.
# 3) Replace the Velocity & Completion Trend AI narrative (slide 3)
#    with the "no insight generated" placeholder text.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- Slide 1: update the generated-on date -------------------------
$slide1 = $p.Slides.Item(1)
$subtitleRange = $slide1.Shapes.Item(2).TextFrame.TextRange
$subtitleRange.Runs(1).Text = "Executive Summary - Generated by AI on 2025-10-12"

# --- Slide 2: replace the "Overall Health & Status" narrative ------
$slide2 = $p.Slides.Item(2)
$slide2Range = $slide2.Shapes.Item(2).TextFrame.TextRange
$slide2Paragraphs = @(
    "Based on the overall health assessment from the various metrics, here are key improvement areas that can be implemented immediately:",
    "PLACEHOLDER_EMPTY_PARAGRAPH",
    "### Immediate Improvement Areas for Program Health:",
    "PLACEHOLDER_EMPTY_PARAGRAPH",
    "1.  **Strengthen Sprint Planning and Commitment Reliability:**",
    "    *   **Evidence:** Consistently low sprint completion rates (e.g., SPRINT-7 committed 76 points, completed 41; SPRINT-8 committed 50, completed 19) and highly inconsistent velocity (e.g., SPRINT-9 at 12 points vs. SPRINT-10 at 52 points). This indicates issues with realistic planning or execution.",
    "    *   **Immediate Action:**",
    "        *   **Implement more rigorous story point estimation:** Review past sprint data during planning to improve accuracy. Encourage the team to account for complexity, unknowns, and potential dependencies when estimating.",
    "        *   **Refine sprint commitment discussions:** Ensure the team has a shared understanding of what `"done`" means for each story. Challenge over-commitment and prioritize a realistic, achievable sprint goal.",
    "        *   **Conduct daily stand-ups with a focus on blockers:** Ensure that progress is tracked daily and any impediments to completing committed work are identified and addressed immediately.",
    "PLACEHOLDER_EMPTY_PARAGRAPH",
    "2.  **Shift-Left Quality Assurance to Reduce Defect Leakage:**",
    "    *   **Evidence:** High defect density across sprints (e.g., SPRINT-6 and SPRINT-7 show 19 defects for 14 stories each). Crucially, a significant number of defects are found late in the cycle: 23 in SIT, 21 in UAT, and 8 in Production. This indicates defects are not being caught early enough.",
    "    *   **Immediate Action:**",
    "        *   **Enhance Definition of Done (DoD):** Include robust unit testing, code reviews by peers, and automated static code analysis as mandatory steps before a story is considered complete by developers.",
    "        *   **Increase Developer Testing:** Encourage developers to write more comprehensive unit and integration tests, and conduct thorough self-testing before handing over to QA.",
    "        *   **Early QA Involvement:** Engage QA engineers earlier in the sprint to review requirements and create test plans, potentially even participating in design discussions to identify potential defect sources.",
    "PLACEHOLDER_EMPTY_PARAGRAPH",
    "3.  **Optimize Team Workload Distribution and Capacity Utilization:**",
    "    *   **Evidence:** Overall team underutilization, with no individual exceeding 70% of their assumed capacity. There's also a significant imbalance, with Fiona (68%) and Bob (61%) having higher loads, while Alice (19%) has a very low utilization.",
    "    *   **Immediate Action:**",
    "        *   **Investigate Low Utilization:** Conduct one-on-one discussions with individuals showing very low utilization (e.g., Alice) to understand reasons – are there skill gaps, blockers, or insufficient task assignment?",
    "        *   **Cross-Training Opportunities:** Identify areas for cross-training within the team to enable more flexible task allocation and reduce reliance on specific individuals for certain types of work.",
    "        *   **Transparent Task Allocation:** Use sprint planning and daily stand-ups to openly discuss task assignments and ensure a more balanced distribution of workload, leveraging the higher capacity available across the team.",
    "PLACEHOLDER_EMPTY_PARAGRAPH",
    "4.  **Proactive Management of RAID Items:**",
    "    *   **Evidence:** A considerable number of open RAID items (17 Risks, 12 Dependencies, 9 Issues, 3 Assumptions). While none are currently overdue, these represent potential future roadblocks that can impact sprint velocity and completion if not addressed.",
    "    *   **Immediate Action:**",
    "        *   **Dedicated RAID Review:** Schedule a dedicated, frequent (e.g., weekly) short meeting to review all open RAID items. Assign clear owners and realistic target dates for mitigation/resolution.",
    "        *   **Escalate Blockers Promptly:** Establish a clear escalation path for critical risks, issues, or unresolved dependencies that are impeding sprint progress. Do not let them linger.",
    "        *   **Communicate Dependencies:** Ensure that external dependencies are clearly communicated to relevant stakeholders with agreed-upon timelines to prevent them from becoming sprint blockers."
)
$slide2Range.Text = [string]::Join("`r", $slide2Paragraphs)

# The source deck keeps a handful of truly blank paragraphs (no run)
# between sections; re-blank those paragraph slots now that the
# paragraph list/structure has been created above.
$slide2Range.Paragraphs(2, 1).Text = ""
$slide2Range.Paragraphs(4, 1).Text = ""
$slide2Range.Paragraphs(11, 1).Text = ""
$slide2Range.Paragraphs(18, 1).Text = ""
$slide2Range.Paragraphs(25, 1).Text = ""

# --- Slide 3: replace the velocity-trend narrative ------------------
$slide3 = $p.Slides.Item(3)
$slide3Range = $slide3.Shapes.Item(2).TextFrame.TextRange
$slide3Paragraphs = @(
    "No specific AI insight generated for this metric in current conversation.",
    "No specific AI insight generated for this metric in current conversation."
)
$slide3Range.Text = [string]::Join("`r", $slide3Paragraphs)

